$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values (A2:A4)
$ws.Range("A2").Value = 29.73
$ws.Range("A3").Value = 11.12
$ws.Range("A4").Value = 41.19

# Add new data row
$ws.Range("A5").Value = 32.020000000000003

# Turn on AutoFilter for the header cell
$ws.Range("A1").AutoFilter() | Out-Null

# Excel registers the filter range as the hidden sheet-scoped
# _FilterDatabase defined name
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet!`$A`$1:`$A`$1")
$fdb.Visible = $false

# Move the active selection to C4
$ws.Range("C4").Select() | Out-Null
